$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Feuil1")
$ws2 = $wb.Worksheets.Item("Feuil2")

# ---- Feuil1 (sheet1) ----
# Row 18: turn N18 into a PRODUCT formula, bump P18, clear U18
$ws1.Range("N18").Formula = "=PRODUCT(O18:P18)"
$ws1.Range("P18").Value = 30
$ws1.Range("U18").ClearContents()

# Row 19: turn N19 into a PRODUCT formula, bump P19, clear U19
$ws1.Range("N19").Formula = "=PRODUCT(O19:P19)"
$ws1.Range("P19").Value = 12
$ws1.Range("U19").ClearContents()

# Row 23: extend the SUM range
$ws1.Range("N23").Formula = "=SUM(N18:N22)"

# Row 25: update total
$ws1.Range("N25").Value = 26400

# ---- Feuil2 (sheet2) ----
# Row 9: update mission details (numeric columns)
$ws2.Range("D9").Value = 7
$ws2.Range("F9").Value = 16

# Row 10: fill in the second mission block (numeric columns)
$ws2.Range("D10").Value = 3
$ws2.Range("F10").Value = 8

# Row 11: fill in the third mission block (numeric columns)
$ws2.Range("D11").Value = 2
$ws2.Range("F11").Value = 6

# ---- Text values that must stay literal text (not auto-converted to
# ----  dates/times/numbers by Excel's smart entry) are written through a
# ----  scratch cell formatted as Text, then pasted as values-only so the
# ----  destination cell keeps its original number format/style.
$sheets  = @($ws1, $ws2, $ws2, $ws2, $ws2, $ws2, $ws2, $ws2, $ws2, $ws2, $ws2, $ws2, $ws2, $ws2, $ws2, $ws2, $ws2, $ws2, $ws2, $ws2, $ws2, $ws2, $ws2)
$refs    = @("L38", "N9", "P9", "B9", "M9", "B10", "L10", "M10", "N10", "O10", "P10", "R10", "S10", "T10", "B11", "L11", "M11", "N11", "O11", "P11", "R11", "S11", "T11")
$vals    = @("2021/08/12", "20/02/2020", "13/02/2020", "522", "23:00", "50", "مهمة ادارية", "23:00", "15/08/2021", "06:00", "12/08/2021", "ولاية الشلف", "بسكرة", "مهمة ادارية", "53", "مهمة ادارية", "23:00", "14/08/2021", "06:00", "12/08/2021", "ولاية باتنة", "بسكرة", "مهمة ادارية")

for ($i = 0; $i -lt $refs.Length; $i++) {
    $sht = $sheets[$i]
    $helper = $sht.Range("ZZ1")
    $helper.NumberFormat = "@"
    $helper.Value = $vals[$i]
    $helper.Copy()
    $sht.Range($refs[$i]).PasteSpecial(-4163)
    $helper.Clear()
}
